$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Style row 3 (A3 date style w/ right+center align, B3/C3 left+center align)
# ---------------------------------------------------------------------------
$ws.Range("A3").HorizontalAlignment = -4152   # xlHAlignRight
$ws.Range("A3").VerticalAlignment = -4108     # xlVAlignCenter

$ws.Range("B3").HorizontalAlignment = -4131   # xlHAlignLeft
$ws.Range("B3").VerticalAlignment = -4108     # xlVAlignCenter

$ws.Range("B3").Copy()
$ws.Range("C3").PasteSpecial(-4122)           # xlPasteFormats

# Propagate the same formats down to rows 4 and 5 ahead of the merge.
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)

$ws.Range("B3").Copy()
$ws.Range("B4:B5").PasteSpecial(-4122)
$ws.Range("C4:C5").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Fix up the misaligned D3/E3/F3 values (Goals / time / reflection)
# ---------------------------------------------------------------------------
$ws.Range("D3").Value = "Goals"
$ws.Range("E3").Value = "12 mins (11:11 - 11:23)"
$ws.Range("F3").Value = "Fine, the goal was previously sketched"

# ---------------------------------------------------------------------------
# 3) New rows describing Users / Similar Systems work
#    (write D5 first so new shared strings land in the same order as target)
# ---------------------------------------------------------------------------
$ws.Range("D5").Value = "Similar Systems"

$ws.Range("D4").Value = "Users"
$ws.Range("E4").NumberFormat = "h:mm"
$ws.Range("E4").Value = "18 mins (11:32 - 11:50)"
$ws.Range("F4").Value = "Nice, the definition of possible users was an interesting task."

# Drop the stray formatted-but-empty E5 cell now that E4 carries the style.
$ws.Range("E5").Clear()

# ---------------------------------------------------------------------------
# 4) Merge the A/B/C columns across the 3-row "App Definition" entry
# ---------------------------------------------------------------------------
$ws.Range("A3:A5").Merge()
$ws.Range("B3:B5").Merge()
$ws.Range("C3:C5").Merge()

# ---------------------------------------------------------------------------
# 5) Move the active selection, matching the authored edit
# ---------------------------------------------------------------------------
[void]$ws.Range("F12").Select()

Write-Host "done"
